# Correct file upload functionality
# Append a new data row (row 85) to each of the four worksheets, mirroring
# the most recent reading that had already been logged in row 84.

$wb = $excel.ActiveWorkbook

function Add-Row85 {
    param($ws, [double]$TimeValue, [string]$B, [string]$C, [string]$D, [string]$E, [int]$F, [double]$G, [int]$H, [int]$I)

    $row = 85

    $ws.Cells.Item($row, 1).Value = $TimeValue
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $B
    $ws.Cells.Item($row, 3).Value = $C
    $ws.Cells.Item($row, 4).Value = $D
    $ws.Cells.Item($row, 5).Value = $E

    $ws.Cells.Item($row, 6).Value = $F
    $ws.Cells.Item($row, 7).Value = $G
    $ws.Cells.Item($row, 8).Value = $H
    $ws.Cells.Item($row, 9).Value = $I
}

$ws1 = $wb.Worksheets.Item("DE_LFT_#1")
$g1 = [double]"7.598631275147109e+23"
Add-Row85 $ws1 45871.43392361111 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x34" "0x14" 380 $g1 308 14

$ws2 = $wb.Worksheets.Item("DE_LFT_#2")
$g2 = [double]"5.68432987514711e+23"
Add-Row85 $ws2 45871.43392361111 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x38" "0xe" 380 $g2 312 14

$ws3 = $wb.Worksheets.Item("DE_PLT_#1")
$g3 = [double]"5.68631262647114e+23"
Add-Row85 $ws3 45871.43392361111 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x77" "0x7" 130 $g3 119 7

$ws4 = $wb.Worksheets.Item("DE_PLT_#2")
$g4 = [double]"9.85046333984776e+23"
Add-Row85 $ws4 45871.43392361111 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x75" "0x3" 130 $g4 117 3
